$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.959.25"
$ws.Range("E2").Value = "  +3.60%  "
$ws.Range("D3").Value = "1.922.57"
$ws.Range("E3").Value = "  +3.51%  "
$ws.Range("E4").Value = "  +3.91%  "
$ws.Range("D5").Value = "323.30"
$ws.Range("E5").Value = "  +3.63%  "
$ws.Range("D6").Value = "1.039"
$ws.Range("E6").Value = "  +3.69%  "
$ws.Range("D7").Value = "0.5241"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("D8").Value = "0.3987"
$ws.Range("E8").Value = "  +4.37%  "
$ws.Range("D9").Value = "0.08426"
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("D10").Value = "1.146"
$ws.Range("E10").Value = "  +3.44%  "
$ws.Range("D11").Value = "42.98"
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("D12").Value = "6.356"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").Value = "1.900.17"
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").Value = "20.76"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").Value = "7.374"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "1.042"
$ws.Range("E16").Value = "  +3.90%  "
$ws.Range("D17").Value = "0.00001124"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").Value = "92.33"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").Value = "0.06888"
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").Value = "18.15"
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("D21").Value = "1.038"
$ws.Range("E21").Value = "  +3.66%  "
$ws.Range("D22").Value = "6.151"
$ws.Range("E22").Value = "  +2.45%  "
$ws.Range("D23").Value = "28.971.07"
$ws.Range("E23").Value = "  +3.50%  "
$ws.Range("D24").Value = "11.35"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").Value = "2.300"
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("D26").Value = "2.123.42"
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("D27").Value = "164.22"
$ws.Range("E27").Value = "  +4.38%  "
$ws.Range("D28").Value = "21.18"
$ws.Range("E28").Value = "  +3.74%  "
$ws.Range("D29").Value = "2.479"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "128.64"
$ws.Range("E30").Value = "  +3.35%  "
$ws.Range("D31").Value = "0.1066"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").Value = "1.066"
$ws.Range("E32").Value = "  +4.10%  "
$ws.Range("D33").Value = "6.024"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").Value = "3.724"
$ws.Range("E34").Value = "  +3.61%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "9.525"
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.02491"
$ws.Range("E36").Value = "  +3.27%  "
$ws.Range("D37").Value = "0.06708"
$ws.Range("E37").Value = "  +3.35%  "
$ws.Range("D38").Value = "0.2238"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").Value = "0.6617"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("D40").Value = "1.268"
$ws.Range("E40").Value = "  +5.03%  "
$ws.Range("D41").Value = "1.206"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("D42").Value = "5.047"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "11.30"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("D44").Value = "0.6202"
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("D45").Value = "13.28"
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("D46").Value = "3.791"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("D47").Value = "1.315"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("D48").Value = "2.044"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").Value = "1.245"
$ws.Range("E49").Value = "  +2.83%  "
$ws.Range("D50").Value = "124.12"
$ws.Range("E50").Value = "  +2.86%  "
$ws.Range("D51").Value = "0.07013"
$ws.Range("E51").Value = "  +3.10%  "
